# Scheduled data refresh: update market price/profit figures across the
# leve-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1104.1451
$ws.Range("J17").Value = 1114.2833
$ws.Range("L17").Value = 3342.8499
$ws.Range("N17").Value = -3678.8499

$ws.Range("H28").Value = 1177.1666
$ws.Range("I28").Value = 1682.2222
$ws.Range("J28").Value = 672.1111
$ws.Range("K28").Value = 1682.2222
$ws.Range("L28").Value = 672.1111
$ws.Range("M28").Value = -1197.2222
$ws.Range("N28").Value = -1642.1111

$ws.Range("H58").Value = 120402.8
$ws.Range("J58").Value = 34000
$ws.Range("L58").Value = 102000
$ws.Range("N58").Value = -102300

$ws.Range("H96").Value = 780.125
$ws.Range("I96").Value = 541.1667
$ws.Range("K96").Value = 1623.5001
$ws.Range("M96").Value = -250.5001

$ws.Range("H98").Value = 2291.9443
$ws.Range("I98").Value = 2291.9443
$ws.Range("K98").Value = 2291.9443
$ws.Range("M98").Value = -793.9443000000001

$ws.Range("H113").Value = 40747268
$ws.Range("I113").Value = 15875261
$ws.Range("J113").Value = 62510276
$ws.Range("K113").Value = 15875261
$ws.Range("L113").Value = 62510276
$ws.Range("M113").Value = -15872007
$ws.Range("N113").Value = -62516784

$ws.Range("H122").Value = 2291.9443
$ws.Range("I122").Value = 2291.9443
$ws.Range("K122").Value = 6875.8329
$ws.Range("M122").Value = -4425.8329

$ws.Range("H138").Value = 3924.449
$ws.Range("I138").Value = 1171.5714
$ws.Range("J138").Value = 7594.952
$ws.Range("K138").Value = 3514.7142
$ws.Range("L138").Value = 22784.856
$ws.Range("M138").Value = 1625.2858
$ws.Range("N138").Value = -33064.856

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 8310.866
$ws.Range("I45").Value = 2639.4285
$ws.Range("K45").Value = 2639.4285
$ws.Range("M45").Value = -2262.4285

$ws.Range("H74").Value = 14585.44
$ws.Range("I74").Value = 20751.838
$ws.Range("J74").Value = 4524.4736
$ws.Range("K74").Value = 20751.838
$ws.Range("L74").Value = 4524.4736
$ws.Range("M74").Value = -19877.838
$ws.Range("N74").Value = -6272.4736

$ws.Range("H77").Value = 14585.44
$ws.Range("I77").Value = 20751.838
$ws.Range("J77").Value = 4524.4736
$ws.Range("K77").Value = 103759.19
$ws.Range("L77").Value = 22622.368
$ws.Range("M77").Value = -99391.19
$ws.Range("N77").Value = -31358.368

$ws.Range("H119").Value = 90000
$ws.Range("J119").Value = 90000
$ws.Range("L119").Value = 90000
$ws.Range("N119").Value = -99676

$ws.Range("H132").Value = 5063.927
$ws.Range("I132").Value = 2840.1482
$ws.Range("K132").Value = 8520.444600000001
$ws.Range("M132").Value = -5990.444600000001

$ws.Range("H134").Value = 40000
$ws.Range("J134").Value = 40000
$ws.Range("L134").Value = 40000
$ws.Range("N134").Value = -50140

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5685.4346
$ws.Range("I134").Value = 2081.739
$ws.Range("J134").Value = 9289.130999999999
$ws.Range("K134").Value = 6245.217000000001
$ws.Range("L134").Value = 27867.393
$ws.Range("M134").Value = -3710.217000000001
$ws.Range("N134").Value = -32937.393

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 7357858.5
$ws.Range("I62").Value = 10421957
$ws.Range("J62").Value = 4021.2
$ws.Range("K62").Value = 10421957
$ws.Range("L62").Value = 4021.2
$ws.Range("M62").Value = -10421333
$ws.Range("N62").Value = -5269.2

$ws.Range("H65").Value = 7357858.5
$ws.Range("I65").Value = 10421957
$ws.Range("J65").Value = 4021.2
$ws.Range("K65").Value = 52109785
$ws.Range("L65").Value = 20106
$ws.Range("M65").Value = -52106665
$ws.Range("N65").Value = -26346

$ws.Range("H86").Value = 5687077.5
$ws.Range("I86").Value = 12505290
$ws.Range("J86").Value = 5233.6665
$ws.Range("K86").Value = 12505290
$ws.Range("L86").Value = 5233.6665
$ws.Range("M86").Value = -12504167
$ws.Range("N86").Value = -7479.6665

$ws.Range("H89").Value = 5687077.5
$ws.Range("I89").Value = 12505290
$ws.Range("J89").Value = 5233.6665
$ws.Range("K89").Value = 62526450
$ws.Range("L89").Value = 26168.3325
$ws.Range("M89").Value = -62520834
$ws.Range("N89").Value = -37400.3325

$ws.Range("H99").Value = 6252.1304
$ws.Range("I99").Value = 5281.5884
$ws.Range("J99").Value = 9002
$ws.Range("K99").Value = 5281.5884
$ws.Range("L99").Value = 9002
$ws.Range("M99").Value = -3783.5884
$ws.Range("N99").Value = -11998

$ws.Range("H126").Value = 6252.1304
$ws.Range("I126").Value = 5281.5884
$ws.Range("J126").Value = 9002
$ws.Range("K126").Value = 15844.7652
$ws.Range("L126").Value = 27006
$ws.Range("M126").Value = -13374.7652
$ws.Range("N126").Value = -31946

$ws.Range("H132").Value = 8085.6523
$ws.Range("I132").Value = 2951.2222
$ws.Range("K132").Value = 8853.6666
$ws.Range("M132").Value = -6323.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 5586.143
$ws.Range("I34").Value = 23
$ws.Range("J34").Value = 6014.077
$ws.Range("K34").Value = 69
$ws.Range("L34").Value = 18042.231
$ws.Range("M34").Value = 15
$ws.Range("N34").Value = -18210.231

$ws.Range("H38").Value = 56.81818
$ws.Range("I38").Value = 61.5
$ws.Range("J38").Value = 54.142857
$ws.Range("K38").Value = 184.5
$ws.Range("L38").Value = 162.428571
$ws.Range("M38").Value = 162.5
$ws.Range("N38").Value = -856.428571

$ws.Range("H125").Value = 6257.143
$ws.Range("J125").Value = 6633.3335
$ws.Range("L125").Value = 19900.0005
$ws.Range("N125").Value = -29740.0005

$ws.Range("H137").Value = 54341.684
$ws.Range("J137").Value = 127339.375
$ws.Range("L137").Value = 382018.125
$ws.Range("N137").Value = -392218.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 21000
$ws.Range("I40").Value = 24000
$ws.Range("K40").Value = 24000
$ws.Range("M40").Value = -23849

$ws.Range("H121").Value = 43165.332
$ws.Range("J121").Value = 43165.332
$ws.Range("L121").Value = 43165.332
$ws.Range("N121").Value = -46659.332

$ws.Range("H122").Value = 4204514
$ws.Range("J122").Value = 4427.75
$ws.Range("L122").Value = 13283.25
$ws.Range("N122").Value = -18183.25

$ws.Range("H126").Value = 4256
$ws.Range("I126").Value = 4193.643
$ws.Range("J126").Value = 4430.6
$ws.Range("K126").Value = 12580.929
$ws.Range("L126").Value = 13291.8
$ws.Range("M126").Value = -10110.929
$ws.Range("N126").Value = -18231.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 422.64285
$ws.Range("I16").Value = 333.63635
$ws.Range("J16").Value = 749
$ws.Range("K16").Value = 333.63635
$ws.Range("L16").Value = 749
$ws.Range("M16").Value = -163.63635
$ws.Range("N16").Value = -1089

$ws.Range("H40").Value = 4930.9165
$ws.Range("I40").Value = 4024.0715
$ws.Range("K40").Value = 4024.0715
$ws.Range("M40").Value = -3888.0715

$ws.Range("H61").Value = 3319.4866
$ws.Range("I61").Value = 1538.08
$ws.Range("J61").Value = 7030.75
$ws.Range("K61").Value = 1538.08
$ws.Range("L61").Value = 7030.75
$ws.Range("M61").Value = -1336.08
$ws.Range("N61").Value = -7434.75

$ws.Range("H75").Value = 43000
$ws.Range("J75").Value = 43000
$ws.Range("L75").Value = 43000
$ws.Range("N75").Value = -44872

$ws.Range("H78").Value = 43000
$ws.Range("J78").Value = 43000
$ws.Range("L78").Value = 129000
$ws.Range("N78").Value = -138360

$ws.Range("H87").Value = 69794.5
$ws.Range("J87").Value = 69794.5
$ws.Range("L87").Value = 69794.5
$ws.Range("N87").Value = -72040.5

$ws.Range("H90").Value = 69794.5
$ws.Range("J90").Value = 69794.5
$ws.Range("L90").Value = 209383.5
$ws.Range("N90").Value = -220615.5

$ws.Range("H113").Value = 3319.4866
$ws.Range("I113").Value = 1538.08
$ws.Range("J113").Value = 7030.75
$ws.Range("K113").Value = 1538.08
$ws.Range("L113").Value = 7030.75
$ws.Range("M113").Value = 631.9200000000001
$ws.Range("N113").Value = -11370.75

$ws.Range("H119").Value = 56900
$ws.Range("J119").Value = 56900
$ws.Range("L119").Value = 56900
$ws.Range("N119").Value = -66576

$ws.Range("H122").Value = 4089.3872
$ws.Range("I122").Value = 2846
$ws.Range("J122").Value = 6700.5
$ws.Range("K122").Value = 8538
$ws.Range("L122").Value = 20101.5
$ws.Range("M122").Value = -6088
$ws.Range("N122").Value = -25001.5

$ws.Range("H132").Value = 12202174
$ws.Range("I132").Value = 25002458
$ws.Range("K132").Value = 75007374
$ws.Range("M132").Value = -75004844

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4138.3125
$ws.Range("I126").Value = 3096.5
$ws.Range("J126").Value = 4763.4
$ws.Range("K126").Value = 9289.5
$ws.Range("L126").Value = 14290.2
$ws.Range("M126").Value = -6819.5
$ws.Range("N126").Value = -19230.2
